# Updates cryptos list data (prices & 1h volume %) per upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (coin names, links, percent-change strings) -- safe to set directly,
# Excel will keep these as text because of non-numeric characters (%, spaces, letters, multiple dots, etc.)
$plainChanges = @{
    'D2' = '27.725.33'
    'E2' = '  +1.00%  '
    'D3' = '1.630.93'
    'E3' = '  -0.33%  '
    'E4' = '  +0.43%  '
    'E6' = '  -1.70%  '
    'E7' = '  +0.51%  '
    'E8' = '  +0.09%  '
    'E9' = '  +1.86%  '
    'E10' = '  +0.24%  '
    'E11' = '  +0.44%  '
    'D12' = '1.862.84'
    'E12' = '  -0.26%  '
    'D13' = '1.648.00'
    'E13' = '  +0.70%  '
    'E14' = '  +0.51%  '
    'E15' = '  -4.32%  '
    'E16' = '  +0.28%  '
    'D17' = '27.764.52'
    'E17' = '  +1.19%  '
    'E18' = '  +0.85%  '
    'D19' = '0.0₃0721'
    'E19' = '  -0.38%  '
    'E20' = '  -0.42%  '
    'E21' = '  +0.31%  '
    'E22' = '  -0.35%  '
    'E23' = '  +1.78%  '
    'E24' = '  +6.44%  '
    'E25' = '  +0.19%  '
    'E27' = '  -1.06%  '
    'E28' = '  +0.34%  '
    'E29' = '  +0.21%  '
    'E30' = '  -0.12%  '
    'E31' = '  -1.20%  '
    'E32' = '  +0.14%  '
    'D33' = '1.466.06'
    'E33' = '  +3.43%  '
    'E34' = '  -2.62%  '
    'E35' = '  -2.66%  '
    'E36' = '  +0.79%  '
    'E37' = '  -0.96%  '
    'E38' = '  +0.08%  '
    'B39' = 'ARBITRUM'
    'C39' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'E39' = '  -0.34%  '
    'B40' = 'TrustWalletToken'
    'C40' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'E40' = '  +1.39%  '
    'E41' = '  +6.61%  '
    'E42' = '  -0.80%  '
    'E43' = '  +0.33%  '
    'E44' = '  -1.31%  '
    'E45' = '  -0.07%  '
    'E46' = '  -1.96%  '
    'D47' = '1.773.32'
    'E47' = '  -0.24%  '
    'E49' = '  +0.23%  '
    'B50' = 'Algorand'
    'C50' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'E50' = '  +0.04%  '
    'B51' = 'EnergySwap'
    'C51' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E51' = '  +1.05%  '
}

foreach ($cellRef in $plainChanges.Keys) {
    $ws.Range($cellRef).Value = $plainChanges[$cellRef]
}

# Price updates that LOOK like plain numbers (e.g. "212.51") -- Excel would normally
# auto-convert these to numeric cells (and mangle the exact decimal text via float
# rounding). Force them to stay text, matching the original inline-string cells, by
# entering them with a leading apostrophe (Excel's standard 'treat as text' marker)
# and then resetting the cell style so no extra formatting/quote-prefix style sticks.
$numericLookingChanges = @{
    'D5' = '212.51'
    'D6' = '0.524'
    'D7' = '1.01'
    'D8' = '22.94'
    'D11' = '0.0890'
    'D14' = '4.04'
    'D16' = '64.38'
    'D18' = '230.70'
    'D20' = '7.59'
    'D23' = '9.94'
    'D25' = '149.83'
    'D26' = '6.91'
    'D27' = '0.111'
    'D29' = '15.58'
    'D31' = '0.0482'
    'D36' = '2.36'
    'D37' = '0.566'
    'D39' = '0.870'
    'D40' = '0.912'
    'D41' = '69.25'
    'D46' = '5.41'
    'D49' = '85.87'
    'D50' = '0.0989'
    'D51' = '7.81'
}

foreach ($cellRef in $numericLookingChanges.Keys) {
    $ws.Range($cellRef).Value = "'" + $numericLookingChanges[$cellRef]
    $ws.Range($cellRef).Style = "Normal"
}
